$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 314
$ws1.Range("F4").Value = 2967
$ws1.Range("F7").Value = 2318
$ws1.Range("F8").Value = 1670
$ws1.Range("F10").Value = 851
$ws1.Range("F13").Value = 2659
$ws1.Range("F15").Value = 1516
$ws1.Range("F16").Value = 7046
$ws1.Range("F18").Value = 7204
$ws1.Range("F19").Value = 2
$ws1.Range("F21").Value = 5437
$ws1.Range("F22").Value = 3107
$ws1.Range("F23").Value = 3476
$ws1.Range("F24").Value = 233
$ws1.Range("F26").Value = 1876
$ws1.Range("F27").Value = 81
$ws1.Range("F30").Value = 23
$ws1.Range("F31").Value = 280
$ws1.Range("F32").Value = 39
$ws1.Range("F33").Value = 2415
$ws1.Range("F34").Value = 1177
$ws1.Range("F35").Value = 2685
$ws1.Range("F36").Value = 23
$ws1.Range("F40").Value = 1068
$ws1.Range("G12").Value = 108

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 215
$ws2.Range("F12").Value = 22

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 314
$ws4.Range("F5").Value = 2967
$ws4.Range("F7").Value = 2318
$ws4.Range("F8").Value = 1670
$ws4.Range("F11").Value = 851
$ws4.Range("F14").Value = 2659
$ws4.Range("F15").Value = 1516
$ws4.Range("F16").Value = 215
$ws4.Range("F19").Value = 7046
$ws4.Range("F21").Value = 7204
$ws4.Range("F23").Value = 5437
$ws4.Range("F24").Value = 3107
$ws4.Range("F25").Value = 3476
$ws4.Range("F27").Value = 233
$ws4.Range("F29").Value = 1876
$ws4.Range("F34").Value = 23
$ws4.Range("F35").Value = 280
$ws4.Range("F36").Value = 39
$ws4.Range("F37").Value = 2415
$ws4.Range("F38").Value = 1177
$ws4.Range("F40").Value = 2685
$ws4.Range("F41").Value = 24
$ws4.Range("F46").Value = 1068
